$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-15 Wednesday" "2024-05-16 Thursday"

Replace-Text "694÷8=86, 6" "286÷4=71, 2"
Replace-Text "269÷7=38, 3" "828÷5=165, 3"
Replace-Text "537÷6=89, 3" "424÷4=106, 0"
Replace-Text "930÷8=116, 2" "444÷8=55, 4"
Replace-Text "900÷5=180, 0" "263÷6=43, 5"

Replace-Text "872÷4=218, 0" "230÷7=32, 6"
Replace-Text "672÷2=336, 0" "171÷9=19, 0"
Replace-Text "453÷3=151, 0" "559÷6=93, 1"
Replace-Text "629÷6=104, 5" "725÷3=241, 2"
Replace-Text "959÷6=159, 5" "291÷3=97, 0"

Replace-Text "640÷5=128, 0" "459÷7=65, 4"
Replace-Text "445÷4=111, 1" "871÷8=108, 7"
Replace-Text "879÷4=219, 3" "813÷4=203, 1"
Replace-Text "215÷6=35, 5" "905÷9=100, 5"
Replace-Text "503÷3=167, 2" "981÷2=490, 1"

Replace-Text "482÷8=60, 2" "939÷3=313, 0"
Replace-Text "167÷2=83, 1" "490÷3=163, 1"
Replace-Text "440÷9=48, 8" "589÷4=147, 1"
Replace-Text "191÷8=23, 7" "370÷3=123, 1"
Replace-Text "521÷6=86, 5" "162÷6=27, 0"

Replace-Text "208÷7=29, 5" "351÷7=50, 1"
Replace-Text "270÷4=67, 2" "666÷5=133, 1"
Replace-Text "109÷8=13, 5" "829÷8=103, 5"
Replace-Text "292÷7=41, 5" "573÷2=286, 1"
Replace-Text "578÷8=72, 2" "323÷4=80, 3"
